# Insert 4 new "Bruno Díaz - Investigación" rows before the existing
# "Federico Speroni" rows 62-64 (shifting them down to 66-68), and append
# one more new row after them (row 69), replacing the previously-empty
# styled D65 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for four new rows above the old row 62 (old rows 62-64 shift
# down to 66-68).
$ws.Rows.Item(62).Resize(4).EntireRow.Insert()

# Carry over the date style (numFmtId 14, "m/d/yyyy") used by the other
# date cells in column B onto the newly inserted cells so no new style
# gets created.
$ws.Range("B61").Copy()
$ws.Range("B62:B65").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 62
$ws.Range("A62").Value = "Bruno Díaz"
$ws.Range("B62").Value = 42889
$ws.Range("C62").Value = 4
$ws.Range("D62").Value = "Sprint 3 - Investigación"
$ws.Range("E62").Value = "Investigación sobre autenticación y lectura de documentación de proyectos (Sprints 3 y 4)"

# New row 63
$ws.Range("A63").Value = "Bruno Díaz"
$ws.Range("B63").Value = 42889
$ws.Range("C63").Value = 3
$ws.Range("D63").Value = "Sprint 3 - Investigación"
$ws.Range("E63").Value = "Investigación sobre animación de carga de datos en cada servicio"

# New row 64
$ws.Range("A64").Value = "Bruno Díaz"
$ws.Range("B64").Value = 42890
$ws.Range("C64").Value = 2
$ws.Range("D64").Value = "Sprint 3 - Investigación"
$ws.Range("E64").Value = "Investigación de pantalla de carga inicial"

# New row 65
$ws.Range("A65").Value = "Bruno Díaz"
$ws.Range("B65").Value = 42890
$ws.Range("C65").Value = 3
$ws.Range("D65").Value = "Sprint 3 - Investigación"
$ws.Range("E65").Value = "Investigación sobre autenticación"

# Rows 66-68 now hold what used to be rows 62-64 (Federico Speroni entries)
# and are already correct after the insert/shift - nothing further to do
# there.

# Append a brand-new row 69 (previously the trailing blank, styled D65
# cell). Reuse the same date style for its date cell too.
$ws.Range("B68").Copy()
$ws.Range("B69").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A69").Value = "Bruno Díaz"
$ws.Range("B69").Value = 42895
$ws.Range("C69").Value = 1
$ws.Range("D69").Value = "Sprint 3 - Investigación"
$ws.Range("E69").Value = "Investigación sobre subida de imágenes"

# Clear the leftover "applyFont" style that used to live on D65 (now
# unused since the blank styled cell was replaced by real data).
$ws.Range("D69").ClearFormats()

# Move the active selection to the new first empty row, matching the
# post-edit state (A70 selected).
$ws.Range("A70").Select()
